$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 4.043363094329834
$ws.Range("B1").Value = 2.528892278671265
$ws.Range("C1").Value = 1.810941100120544
$ws.Range("D1").Value = 1.653499126434326
$ws.Range("E1").Value = 1.690486431121826
